$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header labels in row 1 (D1:J1)
$ws.Range("D1").Value = "Default"
$ws.Range("E1").Value = "MCMC"
$ws.Range("F1").Value = "Mini-Splatting"
$ws.Range("G1").Value = "EAGLES"
$ws.Range("H1").Value = "Mip-Splatting"
$ws.Range("I1").Value = "Gaussian-Pro"
$ws.Range("J1").Value = "Geo-Gaussian"

# Update numeric data values (D2:J7, D9:J10)
$ws.Range("D2").Value = 1.224957701347362
$ws.Range("E2").Value = 0.04344177860779775
$ws.Range("F2").Value = 0.009124921513326629
$ws.Range("G2").Value = -0.9357747657846388
$ws.Range("H2").Value = 0.5037227105672092
$ws.Range("I2").Value = -0.176024416869176
$ws.Range("J2").Value = -0.6694473619648188
$ws.Range("D3").Value = 0.08460405217332961
$ws.Range("E3").Value = 1.091316995387062
$ws.Range("F3").Value = 0.588851806290297
$ws.Range("G3").Value = -1.130979522342304
$ws.Range("H3").Value = 0.8466206464111672
$ws.Range("I3").Value = 0.1478801253159097
$ws.Range("J3").Value = -1.628307072458627
$ws.Range("D4").Value = 1.508158915403333
$ws.Range("E4").Value = 3.304051897373304
$ws.Range("F4").Value = -0.1250320509239081
$ws.Range("G4").Value = -0.2005346945619966
$ws.Range("H4").Value = -7.268493809976309
$ws.Range("I4").Value = 0.30435326465208
$ws.Range("J4").Value = 2.477472606445807
$ws.Range("D5").Value = -2.426513498360832
$ws.Range("E5").Value = -4.119895391172083
$ws.Range("F5").Value = -5.308643137371774
$ws.Range("G5").Value = 5.98508056837265
$ws.Range("H5").Value = 4.414793288220697
$ws.Range("I5").Value = 4.969436954066382
$ws.Range("J5").Value = -3.514271373098783
$ws.Range("D6").Value = -5.231147646353744
$ws.Range("E6").Value = 10.14639296707296
$ws.Range("F6").Value = 11.99327157191658
$ws.Range("G6").Value = -6.913766524134874
$ws.Range("H6").Value = 2.800127194568744
$ws.Range("I6").Value = -6.085657373506661
$ws.Range("J6").Value = -6.709242611018718
$ws.Range("D7").Value = -1.899611007026507
$ws.Range("E7").Value = 6.334340827854865
$ws.Range("F7").Value = -2.351148186562976
$ws.Range("G7").Value = -1.587518152405999
$ws.Range("H7").Value = 5.380464729201268
$ws.Range("I7").Value = -3.194217674659312
$ws.Range("J7").Value = -2.682347832757456
$ws.Range("D9").Value = -0.1129042907591787
$ws.Range("E9").Value = 1.346946205933439
$ws.Range("F9").Value = 3.049268828132393
$ws.Range("G9").Value = 1.452908118253542
$ws.Range("H9").Value = 2.254115055176811
$ws.Range("I9").Value = 1.004839533414738
$ws.Range("J9").Value = -8.995284166447723
$ws.Range("D10").Value = -0.6004657316353234
$ws.Range("E10").Value = -0.5840946745515707
$ws.Range("F10").Value = -0.3202807595768558
$ws.Range("G10").Value = 1.905476895875601
$ws.Range("H10").Value = 0.5281051897351202
$ws.Range("I10").Value = 1.03334950061445
$ws.Range("J10").Value = -1.962093278318804
